# Auto-generated edit script: adds crime data for 2023-05-17
# Updates column J (year 2023) crime totals across Citywide Totals, By Neighborhood,
# and individual neighborhood sheets, plus a few prior-year (2019/2020/2022) corrections.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 2604
$ws.Range('J3').Value = 2662
$ws.Range('F4').Value = 1883
$ws.Range('G4').Value = 1459
$ws.Range('J4').Value = 607
$ws.Range('J5').Value = 210
$ws.Range('I6').Value = 8966
$ws.Range('J6').Value = 3286
$ws.Range('F7').Value = 24073
$ws.Range('G7').Value = 24683
$ws.Range('I7').Value = 26202
$ws.Range('J7').Value = 9369

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J3').Value = 21
$ws.Range('J6').Value = 32
$ws.Range('J7').Value = 96

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 98
$ws.Range('J5').Value = 12
$ws.Range('J7').Value = 318

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J2').Value = 72
$ws.Range('J7').Value = 339

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J3').Value = 67
$ws.Range('J6').Value = 87
$ws.Range('J7').Value = 243

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 73
$ws.Range('J4').Value = 42
$ws.Range('J5').Value = 24
$ws.Range('J7').Value = 282
$ws.Range('J8').Value = 585
$ws.Range('J11').Value = 134
$ws.Range('J15').Value = 111
$ws.Range('J17').Value = 19
$ws.Range('J19').Value = 302
$ws.Range('J20').Value = 191
$ws.Range('J25').Value = 55
$ws.Range('J29').Value = 545
$ws.Range('J33').Value = 382
$ws.Range('J34').Value = 49
$ws.Range('I36').Value = 357
$ws.Range('J36').Value = 141
$ws.Range('J37').Value = 318
$ws.Range('J42').Value = 365
$ws.Range('J47').Value = 80
$ws.Range('J50').Value = 54
$ws.Range('J52').Value = 236
$ws.Range('J55').Value = 112
$ws.Range('F63').Value = 175
$ws.Range('G63').Value = 217
$ws.Range('J63').Value = 46
$ws.Range('J65').Value = 243
$ws.Range('J67').Value = 339
$ws.Range('J78').Value = 131
$ws.Range('J79').Value = 283
$ws.Range('J83').Value = 225
$ws.Range('I85').Value = 1159
$ws.Range('J85').Value = 436
$ws.Range('J88').Value = 94
$ws.Range('J89').Value = 96
$ws.Range('J91').Value = 106
$ws.Range('J94').Value = 80
$ws.Range('J95').Value = 144
$ws.Range('F101').Value = 24073
$ws.Range('G101').Value = 24683
$ws.Range('I101').Value = 26202
$ws.Range('J101').Value = 9369

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J2').Value = 68
$ws.Range('J7').Value = 225

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J2').Value = 54
$ws.Range('J7').Value = 144

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J2').Value = 102
$ws.Range('J3').Value = 115
$ws.Range('J7').Value = 382

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J2').Value = 161
$ws.Range('J3').Value = 180
$ws.Range('J6').Value = 150
$ws.Range('J7').Value = 545

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J3').Value = 84
$ws.Range('J6').Value = 116
$ws.Range('J7').Value = 302

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J3').Value = 166
$ws.Range('J4').Value = 31
$ws.Range('I6').Value = 303
$ws.Range('I7').Value = 1159
$ws.Range('J7').Value = 436

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 73
$ws.Range('J3').Value = 79
$ws.Range('J6').Value = 187
$ws.Range('J7').Value = 365

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('J3').Value = 45
$ws.Range('J7').Value = 131

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J3').Value = 20
$ws.Range('J7').Value = 112

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 46
$ws.Range('J7').Value = 106

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J6').Value = 78
$ws.Range('J7').Value = 283

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('J2').Value = 19
$ws.Range('J6').Value = 21

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('J2').Value = 62
$ws.Range('J7').Value = 191

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range('J5').Value = 1
$ws.Range('J7').Value = 19

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('J2').Value = 51
$ws.Range('J3').Value = 36
$ws.Range('J4').Value = 3
$ws.Range('I6').Value = 112
$ws.Range('I7').Value = 357
$ws.Range('J7').Value = 141

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 53
$ws.Range('J7').Value = 236

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 12
$ws.Range('J7').Value = 49

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J3').Value = 14
$ws.Range('J4').Value = 5
$ws.Range('J6').Value = 43
$ws.Range('J7').Value = 80

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('J2').Value = 24
$ws.Range('J5').Value = 1
$ws.Range('J7').Value = 55

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('J2').Value = 22
$ws.Range('J7').Value = 80

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('J3').Value = 30
$ws.Range('J7').Value = 111

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('J2').Value = 14
$ws.Range('J3').Value = 17
$ws.Range('J7').Value = 54

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J2').Value = 48
$ws.Range('J7').Value = 134

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J6').Value = 20
$ws.Range('J7').Value = 73

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('J3').Value = 28
$ws.Range('J7').Value = 94

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 184
$ws.Range('J3').Value = 190
$ws.Range('J5').Value = 17
$ws.Range('J6').Value = 167
$ws.Range('J7').Value = 585

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('J6').Value = 15
$ws.Range('J7').Value = 24

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J2').Value = 94
$ws.Range('J3').Value = 81
$ws.Range('J7').Value = 282

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('J6').Value = 14
$ws.Range('J7').Value = 42
